$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K2").Value = "`${inventory}"
$ws.Range("K2").Select()
